$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.157.30"
$ws.Range("E2").Value = "  -2.69%  "

$ws.Range("D3").Value = "1.804.04"
$ws.Range("E3").Value = "  -0.59%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.44%  "

$ws.Range("D5").Value = "317.09"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("E6").Value = "  +0.35%  "

$ws.Range("D7").Value = "0.5297"
$ws.Range("E7").Value = "  -2.12%  "

$ws.Range("D8").Value = "0.3773"
$ws.Range("E8").Value = "  -1.69%  "

$ws.Range("D9").Value = "0.07477"
$ws.Range("E9").Value = "  -1.77%  "

$ws.Range("D10").Value = "42.01"
$ws.Range("E10").Value = "  -2.16%  "

$ws.Range("D11").Value = "1.097"
$ws.Range("E11").Value = "  -3.06%  "

$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.42%  "

$ws.Range("D13").Value = "6.214"
$ws.Range("E13").Value = "  -0.52%  "

$ws.Range("D14").Value = "20.54"
$ws.Range("E14").Value = "  -3.84%  "

$ws.Range("D15").Value = "7.352"
$ws.Range("E15").Value = "  -1.86%  "

$ws.Range("D16").Value = "1.796.18"
$ws.Range("E16").Value = "  -0.80%  "

$ws.Range("D17").Value = "89.59"
$ws.Range("E17").Value = "  -3.16%  "

$ws.Range("D18").Value = "0.00001067"
$ws.Range("E18").Value = "  -0.65%  "

$ws.Range("D19").Value = "0.06505"
$ws.Range("E19").Value = "  +0.71%  "

$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.33%  "

$ws.Range("D21").Value = "17.24"
$ws.Range("E21").Value = "  -1.01%  "

$ws.Range("D22").Value = "5.920"
$ws.Range("E22").Value = "  -1.23%  "

$ws.Range("D23").Value = "28.197.14"
$ws.Range("E23").Value = "  -2.52%  "

$ws.Range("D24").Value = "11.16"
$ws.Range("E24").Value = "  -2.83%  "

$ws.Range("D25").Value = "2.091"
$ws.Range("E25").Value = "  -2.04%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "155.87"
$ws.Range("E26").Value = "  -4.19%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "20.47"
$ws.Range("E27").Value = "  -1.27%  "

$ws.Range("D28").Value = "2.006.83"
$ws.Range("E28").Value = "  -0.74%  "

$ws.Range("E29").Value = "  -4.55%  "

$ws.Range("D30").Value = "122.07"
$ws.Range("E30").Value = "  -2.18%  "

$ws.Range("D31").Value = "1.121"
$ws.Range("E31").Value = "  -2.80%  "

$ws.Range("D32").Value = "0.1096"
$ws.Range("E32").Value = "  +7.85%  "

$ws.Range("D33").Value = "5.585"
$ws.Range("E33").Value = "  -4.12%  "

$ws.Range("D34").Value = "3.620"
$ws.Range("E34").Value = "  -1.57%  "

$ws.Range("D35").Value = "0.07229"
$ws.Range("E35").Value = "  +7.99%  "

$ws.Range("D36").Value = "0.2225"
$ws.Range("E36").Value = "  -4.92%  "

$ws.Range("D37").Value = "0.02295"
$ws.Range("E37").Value = "  -1.90%  "

$ws.Range("E38").Value = "  -1.23%  "

$ws.Range("D39").Value = "8.487"
$ws.Range("E39").Value = "  -2.70%  "

$ws.Range("D40").Value = "0.6163"
$ws.Range("E40").Value = "  -3.84%  "

$ws.Range("D41").Value = "11.12"
$ws.Range("E41").Value = "  -5.00%  "

$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "1.434"
$ws.Range("E42").Value = "  +1.72%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "1.181"
$ws.Range("E43").Value = "  -3.75%  "

$ws.Range("D44").Value = "13.41"
$ws.Range("E44").Value = "  -2.24%  "

$ws.Range("D45").Value = "3.683"
$ws.Range("E45").Value = "  -0.23%  "

$ws.Range("D46").Value = "0.5764"
$ws.Range("E46").Value = "  -4.18%  "

$ws.Range("D47").Value = "125.69"
$ws.Range("E47").Value = "  -0.31%  "

$ws.Range("D48").Value = "1.192"
$ws.Range("E48").Value = "  +1.87%  "

$ws.Range("D49").Value = "1.923"
$ws.Range("E49").Value = "  -4.53%  "

$ws.Range("D50").Value = "0.06821"
$ws.Range("E50").Value = "  -2.38%  "

$ws.Range("D51").Value = "72.12"
$ws.Range("E51").Value = "  -2.13%  "
